# NK changes wait fix and other changes
#
# The "Status" column (B) on Sheet1 is rewritten: almost every row now
# reports "Product added to cart- Checkout pending" instead of the old
# per-row success/failure + PR-number messages, except for two rows that
# keep a (possibly reworded) failure message, and row 37 (A=36) which has
# no status at all.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$defaultMsg = "Product added to cart- Checkout pending"
$row14Msg   = "Failed -Invalid Product exist : "
$row27Msg   = "Failed -Invalid Product exist :StarTech.com Portable USB C Multiport Video Adapter - 4k HDMI or VGA, USB 3.0"

for ($r = 2; $r -le 39; $r++) {
    if ($r -eq 37) {
        # no status value for this row
        continue
    }

    if ($r -eq 14) {
        $ws.Cells.Item($r, 2).Value = $row14Msg
    } elseif ($r -eq 27) {
        $ws.Cells.Item($r, 2).Value = $row27Msg
    } else {
        $ws.Cells.Item($r, 2).Value = $defaultMsg
    }
}
